$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (a new, most-recent fiscal-year column
# of data was added to each of the three statements: Income Statement,
# Balance Sheet, Cash Flow Statement). This shifts existing D:K down to E:L.
$ws.Columns("D:D").Insert()

# Carry over the number formats/styles from the (now-shifted) adjacent
# column E into the newly inserted, blank column D for the data rows
# (7 through 102) so the new column matches the existing look (date format
# for the "Period Ending" rows, number format elsewhere).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 143600
$ws.Range("D9").Value2 = "NA"
$ws.Range("D10").Value2 = "NA"
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 23900
$ws.Range("D17").Value2 = 76500
$ws.Range("D18").Value2 = 67100
$ws.Range("D20").Value2 = 16100
$ws.Range("D21").Value2 = 105500
$ws.Range("E21").Value2 = 93400
$ws.Range("D22").Value2 = 0
$ws.Range("D23").Value2 = 83200
$ws.Range("D24").Value2 = 300
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 82900
$ws.Range("D27").Value2 = 82400
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -16100
$ws.Range("D33").Value2 = 82400
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 82400
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 92000
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 800
$ws.Range("D44").Value2 = 200
$ws.Range("D45").Value2 = 800
$ws.Range("D46").Value2 = 0
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("F47").Value2 = 0
$ws.Range("G47").Value2 = 0
$ws.Range("H47").Value2 = 0
$ws.Range("I47").Value2 = 0
$ws.Range("J47").Value2 = 0
$ws.Range("D48").Value2 = 1190700
$ws.Range("D49").Value2 = 19000
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 33100
$ws.Range("E52").Value2 = 27400
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 1343100
$ws.Range("D57").Value2 = 1000
$ws.Range("D58").Value2 = 0
$ws.Range("D59").Value2 = 25000
$ws.Range("D60").Value2 = 0
$ws.Range("D61").Value2 = 615900
$ws.Range("D62").Value2 = "NA"
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 652000
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 46000
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 691100
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 82400
$ws.Range("D83").Value2 = 22400
$ws.Range("E83").Value2 = 21500
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 80900
$ws.Range("D91").Value2 = -268300
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -247000
$ws.Range("E94").Value2 = -80400
$ws.Range("D96").Value2 = -69500
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 190000
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 23900
$ws.Range("E102").Value2 = 42700

Write-Output "done"
